$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.949.58'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '1.907.43'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '0.7978'
$ws.Range('E5').Value = '  +5.17%  '
$ws.Range('D6').Value = '241.59'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.3156'
$ws.Range('E8').Value = '  +3.53%  '
$ws.Range('D9').Value = '26.23'
$ws.Range('E9').Value = '  +4.20%  '
$ws.Range('D10').Value = '0.06913'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').Value = '0.07981'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.899.36'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '0.7405'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '5.188'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '92.82'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '29.943.71'
$ws.Range('D17').Value = '13.96'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').Value = '5.854'
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('D19').Value = '244.69'
$ws.Range('E19').Value = '  +4.89%  '
$ws.Range('D20').Value = '0.000007735'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = '0.9999'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '2.152.64'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '6.801'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').Value = '167.63'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('D26').Value = '9.194'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').Value = '0.1408'
$ws.Range('E27').Value = '  +9.14%  '
$ws.Range('D28').Value = '18.87'
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').Value = '1.364'
$ws.Range('E30').Value = '  +1.88%  '
$ws.Range('D31').Value = '1.517'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').Value = '4.302'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('D34').Value = '0.05492'
$ws.Range('E34').Value = '  +3.91%  '
$ws.Range('E35').Value = '  +1.20%  '
$ws.Range('D36').Value = '0.7285'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '2.720'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('D39').Value = '2.779'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('D40').Value = '6.152'
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('D41').Value = '0.4412'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('D42').Value = '71.91'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '0.8328'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('D45').Value = '1.869'
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').Value = '100.36'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').Value = '7.503'
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('D48').Value = '9.714'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D49').Value = '988.31'
$ws.Range('E49').Value = '  +8.49%  '
$ws.Range('D50').Value = '2.056.83'
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('D51').Value = '36.11'
$ws.Range('E51').Value = '  +0.47%  '
